# Strassen_math.xlsx — "Wrapper base design works"
#
# The 4x4 base matrix on Sheet2 (E6:H9) is replaced with an anti-diagonal
# matrix of 2s, and the workbook's view state ends up back on Sheet2 with
# cell J10 selected.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- New base matrix E6:H9 (anti-diagonal, all 2s) ---------------------
$ws2.Range("E6").Value = 0
$ws2.Range("F6").Value = 2
$ws2.Range("G6").Value = 0
$ws2.Range("H6").Value = 0

$ws2.Range("E7").Value = 0
$ws2.Range("F7").Value = 0
$ws2.Range("G7").Value = 2
$ws2.Range("H7").Value = 0

$ws2.Range("E8").Value = 0
$ws2.Range("F8").Value = 0
$ws2.Range("G8").Value = 0
$ws2.Range("H8").Value = 2

$ws2.Range("E9").Value = 2
$ws2.Range("F9").Value = 0
$ws2.Range("G9").Value = 0
$ws2.Range("H9").Value = 0

# --- Refresh the legacy (Ctrl+Shift+Enter) array formulas that read the
# matrix, so their cached spill values (the non-anchor cells of each
# array range) are recomputed along with it -----------------------------
$ws2.Range("Q6:R7").FormulaArray   = "=E8:F9+G8:H9"
$ws2.Range("Q14:R15").FormulaArray = "=E6:F7"
$ws2.Range("Q22:R23").FormulaArray = "=G8:H9"
$ws2.Range("Q30:R31").FormulaArray = "=E6:F7+G6:H7"

$excel.CalculateFullRebuild()

# --- Refresh the downstream formulas that reference a non-anchor
# (spill) cell of one of those arrays directly, so their cached values
# pick up the refreshed spill results ------------------------------------
$ws2.Range("Z6").Formula  = "=Q7+R7"
$ws2.Range("Z8").Formula  = "=R7"
$ws2.Range("Z14").Formula = "=Q15+R15"
$ws2.Range("Z16").Formula = "=R15"
$ws2.Range("Z22").Formula = "=Q23+R23"
$ws2.Range("Z24").Formula = "=R23"
$ws2.Range("Z30").Formula = "=Q31+R31"
$ws2.Range("Z32").Formula = "=R31"

$excel.CalculateFullRebuild()

# --- View state: Sheet2 becomes the active tab with J10 selected -------
$null = $ws2.Range("J10").Select()
$ws2.Activate()

$wb.Save()
